$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "SvyID " column header (note trailing space) in F1
$ws.Range("F1").Value = "SvyID "

# Update the EndDate (column E) values for the existing rows
$ws.Range("E2").Value = Get-Date -Year 2023 -Month 6 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("E3").Value = Get-Date -Year 2023 -Month 6 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("E4").Value = Get-Date -Year 2023 -Month 7 -Day 24 -Hour 0 -Minute 0 -Second 0

# Populate new SvyID values in column F
$ws.Range("F2").Value = 333
$ws.Range("F3").Value = 10044
$ws.Range("F4").Value = 109

# Move the active selection to match the saved view state
$ws.Range("L12").Select()
